$wb = $excel.ActiveWorkbook

# ALC!row8
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 1332.8572
$ws.Range("I8").Value = 82.5
$ws.Range("J8").Value = 3000
$ws.Range("K8").Value = 247.5
$ws.Range("L8").Value = 9000
$ws.Range("M8").Value = -108.5
$ws.Range("N8").Value = -9278

# ALC!row10
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 4624.875
$ws.Range("I10").Value = 2000
$ws.Range("J10").Value = 4999.857
$ws.Range("K10").Value = 2000
$ws.Range("L10").Value = 4999.857
$ws.Range("M10").Value = -1707
$ws.Range("N10").Value = -5585.857

# ALC!row43
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 924.5
$ws.Range("I43").Value = 450.25
$ws.Range("J43").Value = 1114.2
$ws.Range("K43").Value = 450.25
$ws.Range("L43").Value = 1114.2
$ws.Range("M43").Value = -381.25
$ws.Range("N43").Value = -1252.2

# ALC!row55
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 111.111115
$ws.Range("I55").Value = 116.666664
$ws.Range("J55").Value = 100
$ws.Range("K55").Value = 116.666664
$ws.Range("L55").Value = 100
$ws.Range("M55").Value = 97.333336
$ws.Range("N55").Value = -528

# ALC!row129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1140.9166
$ws.Range("I129").Value = 627.2857
$ws.Range("J129").Value = 1860
$ws.Range("K129").Value = 1881.8571
$ws.Range("L129").Value = 5580
$ws.Range("M129").Value = 3118.1429
$ws.Range("N129").Value = -15580

# ARM!row2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2675897.8
$ws.Range("I2").Value = 2390
$ws.Range("J2").Value = 14706682
$ws.Range("K2").Value = 2390
$ws.Range("L2").Value = 14706682
$ws.Range("M2").Value = -2277
$ws.Range("N2").Value = -14706908

# ARM!row19
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("M19").ClearContents()

# ARM!row45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 51650.55
$ws.Range("I45").Value = 84683.086
$ws.Range("J45").Value = 2101.75
$ws.Range("K45").Value = 84683.086
$ws.Range("L45").Value = 2101.75
$ws.Range("M45").Value = -84306.086
$ws.Range("N45").Value = -2855.75

# ARM!row74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 5006671
$ws.Range("I74").Value = 8333949.5
$ws.Range("J74").Value = 15752.875
$ws.Range("K74").Value = 8333949.5
$ws.Range("L74").Value = 15752.875
$ws.Range("M74").Value = -8333075.5
$ws.Range("N74").Value = -17500.875

# ARM!row77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 5006671
$ws.Range("I77").Value = 8333949.5
$ws.Range("J77").Value = 15752.875
$ws.Range("K77").Value = 41669747.5
$ws.Range("L77").Value = 78764.375
$ws.Range("M77").Value = -41665379.5
$ws.Range("N77").Value = -87500.375

# ARM!row116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 2675897.8
$ws.Range("I116").Value = 2390
$ws.Range("J116").Value = 14706682
$ws.Range("K116").Value = 2390
$ws.Range("L116").Value = 14706682
$ws.Range("M116").Value = -96
$ws.Range("N116").Value = -14711270

# ARM!row122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2288.7646
$ws.Range("I122").Value = 1776.25
$ws.Range("J122").Value = 2744.3333
$ws.Range("K122").Value = 5328.75
$ws.Range("L122").Value = 8232.999899999999
$ws.Range("M122").Value = -2878.75
$ws.Range("N122").Value = -13132.9999

# BSM!row3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2675897.8
$ws.Range("I3").Value = 2390
$ws.Range("J3").Value = 14706682
$ws.Range("K3").Value = 2390
$ws.Range("L3").Value = 14706682
$ws.Range("M3").Value = -2276
$ws.Range("N3").Value = -14706910

# CRP!row62
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3398.75
$ws.Range("I62").Value = 2500
$ws.Range("J62").Value = 3698.3333
$ws.Range("K62").Value = 2500
$ws.Range("L62").Value = 3698.3333
$ws.Range("M62").Value = -1876
$ws.Range("N62").Value = -4946.3333

# CRP!row65
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 3398.75
$ws.Range("I65").Value = 2500
$ws.Range("J65").Value = 3698.3333
$ws.Range("K65").Value = 12500
$ws.Range("L65").Value = 18491.6665
$ws.Range("M65").Value = -9380
$ws.Range("N65").Value = -24731.6665

# CUL!row6
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 100
$ws.Range("I6").Value = 100
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 300
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -187

# CUL!row9
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 900
$ws.Range("I9").Value = 200
$ws.Range("J9").Value = 1250
$ws.Range("K9").Value = 600
$ws.Range("L9").Value = 3750
$ws.Range("M9").Value = -376
$ws.Range("N9").Value = -4198

# CUL!row13
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 649.3333
$ws.Range("I13").Value = 724
$ws.Range("J13").Value = 500
$ws.Range("K13").Value = 2172
$ws.Range("L13").Value = 1500
$ws.Range("M13").Value = -2004
$ws.Range("N13").Value = -1836

# CUL!row16
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 16667
$ws.Range("I16").Value = 20000.5
$ws.Range("J16").Value = 10000
$ws.Range("K16").Value = 60001.5
$ws.Range("L16").Value = 30000
$ws.Range("M16").Value = -59828.5
$ws.Range("N16").Value = -30346

# CUL!row22
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 1000
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 3000
$ws.Range("N22").Value = -3338

# CUL!row27
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H27").Value = 1000
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 1000
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 3000
$ws.Range("N27").Value = -3204

# CUL!row39
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 3990
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 3990
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 11970
$ws.Range("N39").Value = -12558

# CUL!row70
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 4447444
$ws.Range("I70").Value = 5927592
$ws.Range("J70").Value = 7000
$ws.Range("K70").Value = 17782776
$ws.Range("L70").Value = 21000
$ws.Range("M70").Value = -17782461
$ws.Range("N70").Value = -21630

# CUL!row73
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H73").Value = 4447444
$ws.Range("I73").Value = 5927592
$ws.Range("J73").Value = 7000
$ws.Range("K73").Value = 17782776
$ws.Range("L73").Value = 21000
$ws.Range("M73").Value = -17781684
$ws.Range("N73").Value = -23184

# CUL!row100
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H100").Value = 1980
$ws.Range("I100").Value = 1980
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 5940
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -5129
$ws.Range("N100").ClearContents()

# CUL!row113
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 518.73334
$ws.Range("I113").Value = 635.75
$ws.Range("J113").Value = 476.18182
$ws.Range("K113").Value = 1907.25
$ws.Range("L113").Value = 1428.54546
$ws.Range("M113").Value = 262.75
$ws.Range("N113").Value = -5768.54546

# CUL!row122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 14511534
$ws.Range("I122").Value = 25641692
$ws.Range("J122").Value = 42328.8
$ws.Range("K122").Value = 230775228
$ws.Range("L122").Value = 380959.2
$ws.Range("M122").Value = -230772778
$ws.Range("N122").Value = -385859.2

# CUL!row131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1065395
$ws.Range("I131").Value = 2782.5
$ws.Range("J131").Value = 1112622.2
$ws.Range("K131").Value = 8347.5
$ws.Range("L131").Value = 3337866.6
$ws.Range("M131").Value = -3307.5
$ws.Range("N131").Value = -3347946.6

# CUL!row139
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 14708483
$ws.Range("I139").Value = 27778968
$ws.Range("J139").Value = 4187.5
$ws.Range("K139").Value = 83336904
$ws.Range("L139").Value = 12562.5
$ws.Range("M139").Value = -83331764
$ws.Range("N139").Value = -22842.5

# GSM!row13
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").ClearContents()

# GSM!row70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5672.727
$ws.Range("I70").Value = 6480
$ws.Range("J70").Value = 5000
$ws.Range("K70").Value = 6480
$ws.Range("L70").Value = 5000
$ws.Range("M70").Value = -6210
$ws.Range("N70").Value = -5540

# GSM!row73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 5672.727
$ws.Range("I73").Value = 6480
$ws.Range("J73").Value = 5000
$ws.Range("K73").Value = 6480
$ws.Range("L73").Value = 5000
$ws.Range("M73").Value = -5544
$ws.Range("N73").Value = -6872

# LTW!row12
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("N12").ClearContents()

# LTW!row16
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3589.111
$ws.Range("I16").Value = 3589.111
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 3589.111
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -3419.111

# LTW!row68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3092.2307
$ws.Range("I68").Value = 2587.375
$ws.Range("J68").Value = 3900
$ws.Range("K68").Value = 2587.375
$ws.Range("L68").Value = 3900
$ws.Range("M68").Value = -1838.375
$ws.Range("N68").Value = -5398

# LTW!row71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 3092.2307
$ws.Range("I71").Value = 2587.375
$ws.Range("J71").Value = 3900
$ws.Range("K71").Value = 12936.875
$ws.Range("L71").Value = 19500
$ws.Range("M71").Value = -9192.875
$ws.Range("N71").Value = -26988

# LTW!row121
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H121").Value = 0
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()

# LTW!row127
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H127").Value = 54980
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 54980
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 54980
$ws.Range("N127").Value = -64900

# WVR!row10
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H10").Value = 38288.715
$ws.Range("I10").Value = 20003.75
$ws.Range("J10").Value = 62668.668
$ws.Range("K10").Value = 20003.75
$ws.Range("L10").Value = 62668.668
$ws.Range("M10").Value = -19834.75
$ws.Range("N10").Value = -63006.668
